$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 31310.445
$ws.Range("I40").Value = 32724
$ws.Range("J40").Value = 30179.6
$ws.Range("K40").Value = 32724
$ws.Range("L40").Value = 30179.6
$ws.Range("M40").Value = -32549
$ws.Range("N40").Value = -30529.6

$ws.Range("H64").Value = 1453996.1
$ws.Range("J64").Value = 4995.769
$ws.Range("L64").Value = 4995.769
$ws.Range("N64").Value = -5491.769

$ws.Range("H67").Value = 1453996.1
$ws.Range("J67").Value = 4995.769
$ws.Range("L67").Value = 4995.769
$ws.Range("N67").Value = -6711.769

$ws.Range("H107").Value = 421.9
$ws.Range("I107").Value = 439.2857
$ws.Range("K107").Value = 439.2857
$ws.Range("M107").Value = 1480.7143

$ws.Range("H116").Value = 2200720.8
$ws.Range("I116").Value = 3797144
$ws.Range("J116").Value = 5638.5
$ws.Range("K116").Value = 3797144
$ws.Range("L116").Value = 5638.5
$ws.Range("M116").Value = -3793702
$ws.Range("N116").Value = -12522.5

$ws.Range("H131").Value = 4360.647
$ws.Range("I131").Value = 857.7273
$ws.Range("J131").Value = 10782.667
$ws.Range("K131").Value = 2573.1819
$ws.Range("L131").Value = 32348.001
$ws.Range("M131").Value = 2466.8181
$ws.Range("N131").Value = -42428.001

$ws.Range("H132").Value = 106797.21
$ws.Range("I132").Value = 248853.61
$ws.Range("K132").Value = 746560.83
$ws.Range("M132").Value = -744030.83

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -184

$ws.Range("H63").Value = 3047.125
$ws.Range("I63").Value = 2955.4
$ws.Range("K63").Value = 2955.4
$ws.Range("M63").Value = -2269.4

$ws.Range("H66").Value = 3047.125
$ws.Range("I66").Value = 2955.4
$ws.Range("K66").Value = 14777
$ws.Range("M66").Value = -11345

$ws.Range("H109").Value = 128897.5
$ws.Range("J109").Value = 128897.5
$ws.Range("L109").Value = 128897.5
$ws.Range("N109").Value = -131671.5

$ws.Range("H132").Value = 4734.8647
$ws.Range("I132").Value = 1931.8572
$ws.Range("K132").Value = 5795.571599999999
$ws.Range("M132").Value = -3265.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3438.6667
$ws.Range("I20").Value = 3338.3333
$ws.Range("J20").Value = 3639.3333
$ws.Range("K20").Value = 3338.3333
$ws.Range("L20").Value = 3639.3333
$ws.Range("M20").Value = -3091.3333
$ws.Range("N20").Value = -4133.3333

$ws.Range("H82").Value = 9330
$ws.Range("I82").Value = 9330
$ws.Range("K82").Value = 9330
$ws.Range("M82").Value = -8947

$ws.Range("H85").Value = 9330
$ws.Range("I85").Value = 9330
$ws.Range("K85").Value = 9330
$ws.Range("M85").Value = -8004

$ws.Range("H140").Value = 158584.5
$ws.Range("J140").Value = 158584.5
$ws.Range("L140").Value = 158584.5
$ws.Range("N140").Value = -168944.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2106
$ws.Range("I16").Value = 2106
$ws.Range("K16").Value = 2106
$ws.Range("M16").Value = -1819

$ws.Range("H62").Value = 53400.832
$ws.Range("I62").Value = 4999
$ws.Range("J62").Value = 63081.2
$ws.Range("K62").Value = 4999
$ws.Range("L62").Value = 63081.2
$ws.Range("M62").Value = -4375
$ws.Range("N62").Value = -64329.2

$ws.Range("H65").Value = 53400.832
$ws.Range("I65").Value = 4999
$ws.Range("J65").Value = 63081.2
$ws.Range("K65").Value = 24995
$ws.Range("L65").Value = 315406
$ws.Range("M65").Value = -21875
$ws.Range("N65").Value = -321646

$ws.Range("H105").Value = 1421388
$ws.Range("J105").Value = 1999.6666
$ws.Range("L105").Value = 1999.6666
$ws.Range("N105").Value = -5493.6666

$ws.Range("H107").Value = 1013901.5
$ws.Range("I107").Value = 1517198.9
$ws.Range("J107").Value = 7306.6665
$ws.Range("K107").Value = 1517198.9
$ws.Range("L107").Value = 7306.6665
$ws.Range("M107").Value = -1515278.9
$ws.Range("N107").Value = -11146.6665

$ws.Range("H109").Value = 63804.832
$ws.Range("J109").Value = 63804.832
$ws.Range("L109").Value = 63804.832
$ws.Range("N109").Value = -65884.83199999999

$ws.Range("H113").Value = 2106
$ws.Range("I113").Value = 2106
$ws.Range("K113").Value = 2106
$ws.Range("M113").Value = 64

$ws.Range("H134").Value = 2121.2554
$ws.Range("I134").Value = 2150
$ws.Range("K134").Value = 6450
$ws.Range("M134").Value = -3915

$ws.Range("H141").Value = 100292.734
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 100292.734
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 100292.734
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -110652.734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 728.4375
$ws.Range("J22").Value = 1158.375
$ws.Range("L22").Value = 3475.125
$ws.Range("N22").Value = -3813.125

$ws.Range("H27").Value = 728.4375
$ws.Range("J27").Value = 1158.375
$ws.Range("L27").Value = 3475.125
$ws.Range("N27").Value = -3679.125

$ws.Range("H81").Value = 5116.4
$ws.Range("I81").Value = 2791.5
$ws.Range("J81").Value = 6666.3335
$ws.Range("K81").Value = 8374.5
$ws.Range("L81").Value = 19999.0005
$ws.Range("M81").Value = -7251.5
$ws.Range("N81").Value = -22245.0005

$ws.Range("H84").Value = 5116.4
$ws.Range("I84").Value = 2791.5
$ws.Range("J84").Value = 6666.3335
$ws.Range("K84").Value = 25123.5
$ws.Range("L84").Value = 59997.0015
$ws.Range("M84").Value = -19507.5
$ws.Range("N84").Value = -71229.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5295721.5
$ws.Range("J70").Value = 5999.75
$ws.Range("L70").Value = 5999.75
$ws.Range("N70").Value = -6539.75

$ws.Range("H73").Value = 5295721.5
$ws.Range("J73").Value = 5999.75
$ws.Range("L73").Value = 5999.75
$ws.Range("N73").Value = -7871.75

$ws.Range("H123").Value = 53899.75
$ws.Range("J123").Value = 59366.332
$ws.Range("L123").Value = 59366.332
$ws.Range("N123").Value = -64266.332

$ws.Range("H126").Value = 3875.1538
$ws.Range("I126").Value = 2133.3333
$ws.Range("K126").Value = 6399.999899999999
$ws.Range("M126").Value = -3929.999899999999

$ws.Range("H132").Value = 2683
$ws.Range("I132").Value = 2683
$ws.Range("K132").Value = 8049
$ws.Range("M132").Value = -5519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6704.7827
$ws.Range("J46").Value = 7141.263
$ws.Range("L46").Value = 7141.263
$ws.Range("N46").Value = -7517.263

$ws.Range("H132").Value = 3832.3582
$ws.Range("I132").Value = 2839.9792
$ws.Range("K132").Value = 8519.937600000001
$ws.Range("M132").Value = -5989.937600000001

$ws.Range("H136").Value = 3813.4614
$ws.Range("I136").Value = 3297.9167
$ws.Range("K136").Value = 9893.750100000001
$ws.Range("M136").Value = -7343.750100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 731.2273
$ws.Range("I113").Value = 660.61536
$ws.Range("J113").Value = 833.2222
$ws.Range("K113").Value = 1981.84608
$ws.Range("L113").Value = 2499.6666
$ws.Range("M113").Value = 188.15392
$ws.Range("N113").Value = -6839.6666

$ws.Range("H132").Value = 12198358
$ws.Range("I132").Value = 1190.7428
$ws.Range("K132").Value = 3572.2284
$ws.Range("M132").Value = -1042.2284

$ws.Range("H136").Value = 7453.0815
$ws.Range("I136").Value = 3203.08
$ws.Range("J136").Value = 8908.562
$ws.Range("K136").Value = 9609.24
$ws.Range("L136").Value = 26725.686
$ws.Range("M136").Value = -7059.24
$ws.Range("N136").Value = -31825.686
